# Update "想去人数" (F column) values in the "展览" and "全部类型" sheets
# to reflect newly scraped counts.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 19
$ws1.Range("F6").Value = 15
$ws1.Range("F7").Value = 559
$ws1.Range("F8").Value = 7837
$ws1.Range("F9").Value = 748
$ws1.Range("F10").Value = 218
$ws1.Range("F12").Value = 728
$ws1.Range("F13").Value = 27
$ws1.Range("F15").Value = 192
$ws1.Range("F18").Value = 802

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 19
$ws4.Range("F6").Value = 15
$ws4.Range("F8").Value = 559
$ws4.Range("F9").Value = 7837
$ws4.Range("F10").Value = 748
$ws4.Range("F11").Value = 218
$ws4.Range("F13").Value = 728
$ws4.Range("F14").Value = 27
$ws4.Range("F16").Value = 192
$ws4.Range("F19").Value = 802
